# ---------------------------------------------------------------------------
# Hystrix_metadata.xlsx edit:
#   1. Append 3 rows (126-128) to "Column descriptions" describing a new
#      "ll" dataset (columns: species, condylobasal.length, Length.of.the.spine)
#   2. Add a new "Authors" worksheet (after "Reference figure") listing the
#      paper's authors (Surname / First name(s) / ORCID / Affiliation)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Column descriptions" sheet - new dataset "ll"
# ---------------------------------------------------------------------------
$colDesc = $wb.Worksheets.Item("Column descriptions")

$colDesc.Range("A126").Value = "Hystrix.RData"
$colDesc.Range("B126").Value = "ll"
$colDesc.Range("C126").Value = "species"
$colDesc.Range("D126").Value = "Species (they are all from genus Hystrix)"

$colDesc.Range("A127").Value = "Hystrix.RData"
$colDesc.Range("B127").Value = "ll"
$colDesc.Range("C127").Value = "condylobasal.length"
$colDesc.Range("D127").Value = "See Reference figure (sheet 3). Units: cm"
# D127 re-uses the same "See Reference figure..." wording already present
# elsewhere in the sheet (e.g. D13), which carries a slightly different
# font (style index 2). Copy that formatting across so the shared string
# keeps reusing the same look-and-feel.
$colDesc.Range("D13").Copy()
$colDesc.Range("D127").PasteSpecial(-4122) | Out-Null

$colDesc.Range("A128").Value = "Hystrix.RData"
$colDesc.Range("B128").Value = "ll"
$colDesc.Range("C128").Value = "Length.of.the.spine"
$colDesc.Range("D128").Value = "Length of the spinal column. Units: cm"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. New "Authors" worksheet, placed after "Reference figure"
# ---------------------------------------------------------------------------
$refFigure = $wb.Worksheets.Item($wb.Worksheets.Count)
$authors = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $refFigure)
$authors.Name = "Authors"

# Column widths (approximate "best fit" sizing used in the source file)
$authors.Columns.Item(1).ColumnWidth = 10.5
$authors.Columns.Item(2).ColumnWidth = 14.17
$authors.Columns.Item(3).ColumnWidth = 24.83
$authors.Columns.Item(4).ColumnWidth = 140.67

$authorRows = @(
    @("Surname", "First name(s)", "ORCID", "Affiliation"),
    @("Plebani", "Marco", "0000-0001-7064-6550", "School of Life Sciences, University of KwaZulu-Natal, Pietermaritzburg, South Africa "),
    @("Ancillotto", "Leonardo", "NA", "Wildlife Research Unit, Dipartimento di Agraria, Universit a degli Studi di Napoli Federico II, Portici, Italy"),
    @("Lovari", "Sandro", "NA", "Unità di Ricerca in Ecologia Comportamentale, Etologia e Gestione della Fauna, Dipartimento di Scienze della Vita, Università di Siena, Siena, Italy"),
    @("Russo", "Danilo", "0000-0002-1934-7130", "Wildlife Research Unit, Dipartimento di Agraria, Università degli Studi di Napoli Federico II"),
    @("Nerva", "Luca", "0000-0001-5009-5798", "Research Centre for Viticulture and Enology, CREA, Conegliano Veneto (Treviso), Italy & Institute for Sustainable Plant Protection, Torino, Italy"),
    @("Mohamed", "Walid Fathy", "NA", "Department of Biological and Geological Sciences, Faculty of Education, Ain Shams University, Roxy, Cairo, Egypt "),
    @("Motro", "Yoav", "NA", "Vertebrate and Snail Division, Ministry of Agriculture and Rural Development, Rishon Leziyyon, Israel"),
    @("Di Bari", "Pietro", "NA", "Dipartimento di Scienze Biologiche, Geologiche e Ambientali – Università di Catania, Via Androne 81, 95131, Catania, Italia"),
    @("Mori", "Emiliano", "0000-0001-8108-7950", "Consiglio Nazionale delle Ricerche, Istituto di Ricerca sugli Ecosistemi Terrestri, Via Madonna del Piano 10, 50019, Sesto Fiorentino (FI), Italy")
)

for ($i = 0; $i -lt $authorRows.Length; $i++) {
    $r = $i + 1
    $rowData = $authorRows[$i]

    $cellA = $authors.Cells.Item($r, 1)
    $cellA.Value = $rowData[0]
    $cellA.Font.Name = "Calibri (Body)"

    $cellB = $authors.Cells.Item($r, 2)
    $cellB.Value = $rowData[1]
    $cellB.Font.Name = "Calibri (Body)"

    $cellC = $authors.Cells.Item($r, 3)
    $cellC.Value = $rowData[2]
    $cellC.Font.Name = "Calibri (Body)"

    $cellD = $authors.Cells.Item($r, 4)
    $cellD.Value = $rowData[3]
    $cellD.Font.Name = "Calibri (Body)"
}

# Row 2 (Plebani)'s ORCID was entered/wrapped differently from the rest -
# reproduce that distinct look (wrapped text, top aligned, text format).
$orcidPlebani = $authors.Range("C2")
$orcidPlebani.WrapText = $true
$orcidPlebani.VerticalAlignment = -4160
$orcidPlebani.NumberFormat = "@"

# A handful of affiliation / ORCID cells carry an explicit black font color
# (as opposed to the default automatic/theme color used elsewhere).
$blackFontCells = @("D5", "D6", "D9", "C10", "D10")
foreach ($addr in $blackFontCells) {
    $cell = $authors.Range($addr)
    $cell.Font.Color = 0
}

$authors.Range("A10").Select()
$authors.Activate()

$excel.CutCopyMode = 0
